# Auto-generated Excel COM-interop script
# Applies the '2025-10-20' crime data update across 29 worksheets (175 cell edits)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("C2").Value = 59
$ws.Range("B3").Value = 69
$ws.Range("D3").Value = 114
$ws.Range("G3").Value = 122
$ws.Range("H3").Value = 123
$ws.Range("J3").Value = 188
$ws.Range("H7").Value = 5
$ws.Range("G9").Value = 396
$ws.Range("H9").Value = 388
$ws.Range("J9").Value = 356
$ws.Range("K9").Value = 431
$ws.Range("L9").Value = 385
$ws.Range("B10").Value = 1120
$ws.Range("C10").Value = 1334
$ws.Range("D10").Value = 1527
$ws.Range("E10").Value = 1814
$ws.Range("F10").Value = 1842
$ws.Range("H10").Value = 502
$ws.Range("I10").Value = 739
$ws.Range("J10").Value = 606
$ws.Range("K10").Value = 592
$ws.Range("B11").Value = 1559
$ws.Range("C11").Value = 1888
$ws.Range("D11").Value = 2077
$ws.Range("E11").Value = 2382
$ws.Range("F11").Value = 2474
$ws.Range("G11").Value = 1422
$ws.Range("H11").Value = 1124
$ws.Range("I11").Value = 1478
$ws.Range("J11").Value = 1275
$ws.Range("K11").Value = 1364
$ws.Range("L11").Value = 1295

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range("H3").Value = 8
$ws.Range("K7").Value = 23
$ws.Range("H9").Value = 69
$ws.Range("K9").Value = 66

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range("H9").Value = 11
$ws.Range("H10").Value = 21

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range("G3").Value = 7
$ws.Range("G9").Value = 40

$ws = $wb.Worksheets.Item('Loop')
$ws.Range("H6").Value = 2
$ws.Range("G8").Value = 51
$ws.Range("H8").Value = 55
$ws.Range("K8").Value = 42
$ws.Range("B9").Value = 164
$ws.Range("D9").Value = 448
$ws.Range("E9").Value = 534
$ws.Range("H9").Value = 90
$ws.Range("I9").Value = 171
$ws.Range("J9").Value = 97
$ws.Range("B10").Value = 205
$ws.Range("D10").Value = 512
$ws.Range("E10").Value = 602
$ws.Range("G10").Value = 224
$ws.Range("H10").Value = 170
$ws.Range("I10").Value = 285
$ws.Range("J10").Value = 199
$ws.Range("K10").Value = 175

$ws = $wb.Worksheets.Item('Old Town')
$ws.Range("C7").Value = 22
$ws.Range("H7").Value = 10
$ws.Range("C8").Value = 29
$ws.Range("H8").Value = 20

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range("B7").Value = 15
$ws.Range("D7").Value = 29
$ws.Range("B8").Value = 29
$ws.Range("D8").Value = 47

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("C10").Value = 9
$ws.Range("H19").Value = 21
$ws.Range("C28").Value = 118
$ws.Range("K29").Value = 25
$ws.Range("D30").Value = 28
$ws.Range("F30").Value = 11
$ws.Range("H32").Value = 69
$ws.Range("K32").Value = 66
$ws.Range("G36").Value = 40
$ws.Range("D47").Value = 45
$ws.Range("J52").Value = 26
$ws.Range("B53").Value = 205
$ws.Range("D53").Value = 512
$ws.Range("E53").Value = 602
$ws.Range("G53").Value = 224
$ws.Range("H53").Value = 170
$ws.Range("I53").Value = 285
$ws.Range("J53").Value = 199
$ws.Range("K53").Value = 175
$ws.Range("D61").Value = 23
$ws.Range("C62").Value = 21
$ws.Range("D62").Value = 24
$ws.Range("B65").Value = 29
$ws.Range("D65").Value = 47
$ws.Range("I68").Value = 8
$ws.Range("C70").Value = 29
$ws.Range("H70").Value = 20
$ws.Range("B71").Value = 5
$ws.Range("B74").Value = 42
$ws.Range("C74").Value = 36
$ws.Range("D74").Value = 79
$ws.Range("E74").Value = 70
$ws.Range("F74").Value = 83
$ws.Range("L77").Value = 44
$ws.Range("D78").Value = 57
$ws.Range("B81").Value = 27
$ws.Range("E83").Value = 20
$ws.Range("C88").Value = 15
$ws.Range("F89").Value = 15
$ws.Range("J92").Value = 24
$ws.Range("J93").Value = 5
$ws.Range("B95").Value = 14
$ws.Range("E95").Value = 83
$ws.Range("B96").Value = 16
$ws.Range("B99").Value = 1559
$ws.Range("C99").Value = 1888
$ws.Range("D99").Value = 2077
$ws.Range("E99").Value = 2382
$ws.Range("F99").Value = 2474
$ws.Range("G99").Value = 1422
$ws.Range("H99").Value = 1124
$ws.Range("I99").Value = 1478
$ws.Range("J99").Value = 1275
$ws.Range("K99").Value = 1364
$ws.Range("L99").Value = 1295

$ws = $wb.Worksheets.Item('Washington Park')
$ws.Range("F6").Value = 12
$ws.Range("F7").Value = 15

$ws = $wb.Worksheets.Item('South Chicago')
$ws.Range("B5").Value = 13
$ws.Range("B6").Value = 27

$ws = $wb.Worksheets.Item('Rush & Division')
$ws.Range("D5").Value = 53
$ws.Range("D6").Value = 57

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range("C8").Value = 76
$ws.Range("C9").Value = 118

$ws = $wb.Worksheets.Item('Lake View')
$ws.Range("D7").Value = 35
$ws.Range("D8").Value = 45

$ws = $wb.Worksheets.Item('Fuller Park')
$ws.Range("K8").Value = 7
$ws.Range("K9").Value = 25

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range("E5").Value = 12
$ws.Range("E6").Value = 20

$ws = $wb.Worksheets.Item('River North')
$ws.Range("B6").Value = 39
$ws.Range("C6").Value = 30
$ws.Range("D6").Value = 69
$ws.Range("E6").Value = 64
$ws.Range("F6").Value = 73
$ws.Range("B7").Value = 42
$ws.Range("C7").Value = 36
$ws.Range("D7").Value = 79
$ws.Range("E7").Value = 70
$ws.Range("F7").Value = 83

$ws = $wb.Worksheets.Item('West Loop')
$ws.Range("J8").Value = 20
$ws.Range("J9").Value = 24

$ws = $wb.Worksheets.Item('Logan Square')
$ws.Range("J6").Value = 6
$ws.Range("J8").Value = 26

$ws = $wb.Worksheets.Item('Near South Side')
$ws.Range("C7").Value = 18
$ws.Range("D7").Value = 20
$ws.Range("C8").Value = 21
$ws.Range("D8").Value = 24

$ws = $wb.Worksheets.Item('West Pullman')
$ws.Range("J3").Value = 4
$ws.Range("J6").Value = 5

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range("L8").Value = 20
$ws.Range("L10").Value = 44

$ws = $wb.Worksheets.Item('West Town')
$ws.Range("B6").Value = 10
$ws.Range("E6").Value = 74
$ws.Range("B7").Value = 14
$ws.Range("E7").Value = 83

$ws = $wb.Worksheets.Item('Wicker Park')
$ws.Range("B6").Value = 16
$ws.Range("B7").Value = 16

$ws = $wb.Worksheets.Item('O''Hare')
$ws.Range("I7").Value = 4
$ws.Range("I8").Value = 8

$ws = $wb.Worksheets.Item('Avondale')
$ws.Range("C6").Value = 7
$ws.Range("C7").Value = 9

$ws = $wb.Worksheets.Item('Gage Park')
$ws.Range("D3").Value = 1
$ws.Range("F7").Value = 9
$ws.Range("D8").Value = 28
$ws.Range("F8").Value = 11

$ws = $wb.Worksheets.Item('Bucktown')
$ws.Range("B3").Value = 1

$ws = $wb.Worksheets.Item('Portage Park')
$ws.Range("B6").Value = 5

$ws = $wb.Worksheets.Item('Washington Heights')
$ws.Range("C2").Value = 1
$ws.Range("C7").Value = 15
